# Auto-generated Excel COM-interop script to update market price / profit
# columns (H:N) across multiple worksheets, per scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H45").Value = 1012.7778
$ws.Range("I45").Value = 749.5
$ws.Range("J45").Value = 1539.3334
$ws.Range("K45").Value = 2248.5
$ws.Range("L45").Value = 4618.0002
$ws.Range("M45").Value = -2056.5
$ws.Range("N45").Value = -5002.0002
$ws.Range("H62").Value = 6640.3
$ws.Range("I62").Value = 3079.4
$ws.Range("J62").Value = 10201.2
$ws.Range("K62").Value = 3079.4
$ws.Range("L62").Value = 10201.2
$ws.Range("M62").Value = -2455.4
$ws.Range("N62").Value = -11449.2
$ws.Range("H65").Value = 6640.3
$ws.Range("I65").Value = 3079.4
$ws.Range("J65").Value = 10201.2
$ws.Range("K65").Value = 15397
$ws.Range("L65").Value = 51006
$ws.Range("M65").Value = -12277
$ws.Range("N65").Value = -57246
$ws.Range("H98").Value = 335629
$ws.Range("I98").Value = 1954.6666
$ws.Range("K98").Value = 1954.6666
$ws.Range("M98").Value = -456.6666
$ws.Range("H112").Value = 1854.9565
$ws.Range("J112").Value = 1881.091
$ws.Range("L112").Value = 5643.272999999999
$ws.Range("N112").Value = -7859.272999999999
$ws.Range("H122").Value = 335629
$ws.Range("I122").Value = 1954.6666
$ws.Range("K122").Value = 5863.9998
$ws.Range("M122").Value = -3413.9998
$ws.Range("H132").Value = 1465.2667
$ws.Range("I132").Value = 1290.3334
$ws.Range("J132").Value = 2165
$ws.Range("K132").Value = 3871.0002
$ws.Range("L132").Value = 6495
$ws.Range("M132").Value = -1341.0002
$ws.Range("N132").Value = -11555
$ws.Range("H135").Value = 1097.6364
$ws.Range("I135").Value = 907.4
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 8166.599999999999
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -5631.599999999999
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 12050709
$ws.Range("I137").Value = 38463070
$ws.Range("J137").Value = 2964.5088
$ws.Range("K137").Value = 115389210
$ws.Range("L137").Value = 8893.526400000001
$ws.Range("M137").Value = -115386660
$ws.Range("N137").Value = -13993.5264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 6852.05
$ws.Range("I2").Value = 1346.2142
$ws.Range("K2").Value = 1346.2142
$ws.Range("M2").Value = -1233.2142
$ws.Range("H32").Value = 3173.8823
$ws.Range("I32").Value = 2134.825
$ws.Range("J32").Value = 19798.8
$ws.Range("K32").Value = 2134.825
$ws.Range("L32").Value = 19798.8
$ws.Range("M32").Value = -1847.825
$ws.Range("N32").Value = -20372.8
$ws.Range("H61").Value = 3562.4443
$ws.Range("I61").Value = 2531.15
$ws.Range("K61").Value = 2531.15
$ws.Range("M61").Value = -2319.15
$ws.Range("H116").Value = 6852.05
$ws.Range("I116").Value = 1346.2142
$ws.Range("K116").Value = 1346.2142
$ws.Range("M116").Value = 947.7858000000001
$ws.Range("H122").Value = 3749.3823
$ws.Range("I122").Value = 3198.76
$ws.Range("K122").Value = 9596.280000000001
$ws.Range("M122").Value = -7146.280000000001
$ws.Range("H132").Value = 2380.9556
$ws.Range("I132").Value = 1697.2593
$ws.Range("K132").Value = 5091.7779
$ws.Range("M132").Value = -2561.7779
$ws.Range("H135").Value = 46527.75
$ws.Range("J135").Value = 46527.75
$ws.Range("L135").Value = 46527.75
$ws.Range("N135").Value = -56667.75
$ws.Range("H136").Value = 3562.4443
$ws.Range("I136").Value = 2531.15
$ws.Range("K136").Value = 7593.450000000001
$ws.Range("M136").Value = -5043.450000000001
$ws.Range("H137").Value = 71995
$ws.Range("J137").Value = 71995
$ws.Range("L137").Value = 71995
$ws.Range("N137").Value = -82195

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 6852.05
$ws.Range("I3").Value = 1346.2142
$ws.Range("K3").Value = 1346.2142
$ws.Range("M3").Value = -1232.2142
$ws.Range("H33").Value = 20000
$ws.Range("J33").Value = 20000
$ws.Range("L33").Value = 20000
$ws.Range("N33").Value = -20672
$ws.Range("H132").Value = 61299.5
$ws.Range("J132").Value = 61299.5
$ws.Range("L132").Value = 61299.5
$ws.Range("N132").Value = -71419.5
$ws.Range("H134").Value = 3076.7273
$ws.Range("I134").Value = 1647.7778
$ws.Range("J134").Value = 9507
$ws.Range("K134").Value = 4943.3334
$ws.Range("L134").Value = 28521
$ws.Range("M134").Value = -2408.3334
$ws.Range("N134").Value = -33591
$ws.Range("H137").Value = 69748.25
$ws.Range("J137").Value = 69748.25
$ws.Range("L137").Value = 69748.25
$ws.Range("N137").Value = -79948.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 8760
$ws.Range("I6").Value = 8760
$ws.Range("K6").Value = 8760
$ws.Range("M6").Value = -8647
$ws.Range("H12").Value = 1336669
$ws.Range("J12").Value = 5003.5
$ws.Range("L12").Value = 5003.5
$ws.Range("N12").Value = -5343.5
$ws.Range("H31").Value = 32603.441
$ws.Range("I31").Value = 1922.6296
$ws.Range("J31").Value = 150943.72
$ws.Range("K31").Value = 1922.6296
$ws.Range("L31").Value = 150943.72
$ws.Range("M31").Value = -1627.6296
$ws.Range("N31").Value = -151533.72
$ws.Range("H32").Value = 7999.5
$ws.Range("I32").Value = 4999
$ws.Range("J32").Value = 11000
$ws.Range("K32").Value = 4999
$ws.Range("L32").Value = 11000
$ws.Range("M32").Value = -4683
$ws.Range("N32").Value = -11632
$ws.Range("H34").Value = 32603.441
$ws.Range("I34").Value = 1922.6296
$ws.Range("J34").Value = 150943.72
$ws.Range("K34").Value = 1922.6296
$ws.Range("L34").Value = 150943.72
$ws.Range("M34").Value = -1720.6296
$ws.Range("N34").Value = -151347.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 31996408
$ws.Range("I4").Value = 40740910
$ws.Range("J4").Value = 10532633
$ws.Range("K4").Value = 122222730
$ws.Range("L4").Value = 31597899
$ws.Range("M4").Value = -122222618
$ws.Range("N4").Value = -31598123
$ws.Range("H5").Value = 2565490.5
$ws.Range("I5").Value = 504.44446
$ws.Range("K5").Value = 1513.33338
$ws.Range("M5").Value = -1401.33338
$ws.Range("H46").Value = 2164.9583
$ws.Range("I46").Value = 372.66666
$ws.Range("J46").Value = 3240.3333
$ws.Range("K46").Value = 1117.99998
$ws.Range("L46").Value = 9720.999899999999
$ws.Range("M46").Value = -1026.99998
$ws.Range("N46").Value = -9902.999899999999
$ws.Range("H80").Value = 8080.5
$ws.Range("I80").Value = 7128.4
$ws.Range("K80").Value = 21385.2
$ws.Range("M80").Value = -20449.2
$ws.Range("H83").Value = 8080.5
$ws.Range("I83").Value = 7128.4
$ws.Range("K83").Value = 64155.6
$ws.Range("M83").Value = -59475.6
$ws.Range("H121").Value = 1276.2
$ws.Range("J121").Value = 1233.1666
$ws.Range("L121").Value = 3699.4998
$ws.Range("N121").Value = -6319.4998
$ws.Range("H131").Value = 4150314.8
$ws.Range("J131").Value = 2819779.5
$ws.Range("L131").Value = 8459338.5
$ws.Range("N131").Value = -8469418.5
$ws.Range("H132").Value = 4062.2273
$ws.Range("I132").Value = 3249
$ws.Range("K132").Value = 29241
$ws.Range("M132").Value = -26711
$ws.Range("H135").Value = 2565490.5
$ws.Range("I135").Value = 504.44446
$ws.Range("K135").Value = 4540.00014
$ws.Range("M135").Value = -2005.00014
$ws.Range("H140").Value = 2592.6365
$ws.Range("I140").Value = 1750.8889
$ws.Range("J140").Value = 6380.5
$ws.Range("K140").Value = 5252.6667
$ws.Range("L140").Value = 19141.5
$ws.Range("M140").Value = -72.66669999999976
$ws.Range("N140").Value = -29501.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 306.25
$ws.Range("I9").Value = 375
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 375
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -205
$ws.Range("N9").Value = -440
$ws.Range("H132").Value = 2528.4412
$ws.Range("I132").Value = 2046.04
$ws.Range("J132").Value = 3868.4443
$ws.Range("K132").Value = 6138.12
$ws.Range("L132").Value = 11605.3329
$ws.Range("M132").Value = -3608.12
$ws.Range("N132").Value = -16665.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 5000000
$ws.Range("I14").Value = 5000000
$ws.Range("K14").Value = 5000000
$ws.Range("M14").Value = -4999828
$ws.Range("H32").Value = 4000
$ws.Range("I32").Value = 4000
$ws.Range("K32").Value = 4000
$ws.Range("M32").Value = -3683
$ws.Range("H46").Value = 6859.864
$ws.Range("I46").Value = 3666.5
$ws.Range("J46").Value = 9521
$ws.Range("K46").Value = 3666.5
$ws.Range("L46").Value = 9521
$ws.Range("M46").Value = -3478.5
$ws.Range("N46").Value = -9897
$ws.Range("H132").Value = 3962.2173
$ws.Range("I132").Value = 3278.4546
$ws.Range("J132").Value = 19005
$ws.Range("K132").Value = 9835.363799999999
$ws.Range("L132").Value = 57015
$ws.Range("M132").Value = -7305.363799999999
$ws.Range("N132").Value = -62075
$ws.Range("H136").Value = 2780
$ws.Range("I136").Value = 1875.6976
$ws.Range("K136").Value = 5627.0928
$ws.Range("M136").Value = -3077.0928

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 40000
$ws.Range("J44").Value = 40000
$ws.Range("L44").Value = 40000
$ws.Range("N44").Value = -41108
$ws.Range("H126").Value = 4317.5386
$ws.Range("I126").Value = 4510.25
$ws.Range("K126").Value = 13530.75
$ws.Range("M126").Value = -11060.75
$ws.Range("H132").Value = 2120.0286
$ws.Range("I132").Value = 1655.1034
$ws.Range("K132").Value = 4965.3102
$ws.Range("M132").Value = -2435.3102
$ws.Range("H136").Value = 2985.2075
$ws.Range("I136").Value = 2039.4286
$ws.Range("K136").Value = 6118.2858
$ws.Range("M136").Value = -3568.2858
